# Revert "Drop in all data files from 3.0 RMI script"
# Re-introduce the "Texas Notes" sheet (with reviewer commentary) between
# "Calculations" and "PEUDfSbQL", and restore the view/selection state that
# goes with it (PEUDfSbQL becomes the active/selected tab again).

$wb = $excel.ActiveWorkbook

# Insert a new worksheet immediately before "PEUDfSbQL" and rename it.
$peSheet = $wb.Worksheets.Item("PEUDfSbQL")
$notes = $wb.Worksheets.Add($peSheet)
$notes.Name = "Texas Notes"

# Reviewer's notes, one sentence/fragment per cell (blank rows intentional).
$notes.Range("A1").Value = "They are just comparing the efficiency of new appliances:"
$notes.Range("A2").Value = "standard versus energy star rebate qualifying"
$notes.Range("A3").Value = "i.e., the point of this spreadsheet is to estimate how much a household's energy consumption would change"
$notes.Range("A4").Value = "if they decide to use a rebate to get a higher-efficiency appliance instead of just opting for the cheap alternative. "
$notes.Range("A6").Value = "I think the method makes sense"
$notes.Range("A8").Value = "And there's no reason that Texas should be different. "
$notes.Range("A9").Value = "New technology in Texas should be as efficient as new technology across the US. "
$notes.Range("A10").Value = "The only difference might be if Texas rebates incentivize a different level of efficiency than"
$notes.Range("A11").Value = "national rebates do, but some of the other sources used in the building input files seem"
$notes.Range("A12").Value = "to indicate that Texas doesn't usually have appliance rebates on top of the national ones. "
$notes.Range("A14").Value = "So it's a good assumption that if a Texan uses a rebate to buy a more efficient appliance, that"
$notes.Range("A15").Value = "rebate will be a national one and it will be based on national energy star standards."

$notes.Range("A16").Select()

# Restore the per-sheet cursor/selection positions.
$wb.Worksheets.Item("About").Range("E25").Select()
$wb.Worksheets.Item("Data").Range("C6").Select()
$wb.Worksheets.Item("Calculations").Range("B18").Select()

# PEUDfSbQL ends up the active tab, with B7 selected.
$wb.Worksheets.Item("PEUDfSbQL").Activate()
$wb.Worksheets.Item("PEUDfSbQL").Range("B7").Select()
